$d = $word.ActiveDocument

# --- 1. Replace the first figure (Gross Plot Ratio) with a hyperlink to the image URL ---
$shp1 = $d.InlineShapes.Item(1)
$r1 = $shp1.Range
$shp1.Delete()
$d.Hyperlinks.Add($r1, "https://ura.gov.sg/-/media/Corporate/Guidelines/Development-control/Flats-Condominiums/F01_Gross_Plot_Ratio.jpg?h=100%25&w=100%25", $null, $null, "https://ura.gov.sg/-/media/Corporate/Guidelines/Development-control/Flats-Condominiums/F01_Gross_Plot_Ratio.jpg?h=100%25&w=100%25") | Out-Null

# --- 2. Replace the second figure (B2-White site) with a hyperlink to the image URL ---
$shp2 = $d.InlineShapes.Item(1)
$r2 = $shp2.Range
$shp2.Delete()
$d.Hyperlinks.Add($r2, "https://ura.gov.sg/-/media/Corporate/Guidelines/Development-control/Industrial/B2-White-Site.jpg?h=100%25&w=100%25", $null, $null, "https://ura.gov.sg/-/media/Corporate/Guidelines/Development-control/Industrial/B2-White-Site.jpg?h=100%25&w=100%25") | Out-Null
